# "Corrections after week10 class"
# The diagram on slide 1 numbers seven textboxes/shapes around a picture.
# The numbering was reversed (1..7 -> 7..1) as a correction; the textbox
# holding "1" (which becomes "7") also needs to grow slightly to fit the
# wider glyph, matching PowerPoint's auto-fit behaviour for the textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)
$items = $grp.GroupItems

for ($i = 1; $i -le $items.Count; $i++) {
    $sh = $items.Item($i)

    if ($sh.Id -eq 31) {
        # TextBox 30: "1" -> "7" (textbox widens to fit the new digit)
        $sh.Width = 13.374
        $sh.TextFrame.TextRange.Text = "7"
    }
    elseif ($sh.Id -eq 32) {
        # TextBox 31: "3" -> "5"
        $sh.TextFrame.TextRange.Text = "5"
    }
    elseif ($sh.Id -eq 33) {
        # TextBox 32: "2" -> "6"
        $sh.TextFrame.TextRange.Text = "6"
    }
    elseif ($sh.Id -eq 35) {
        # TextBox 34: "5" -> "3"
        $sh.TextFrame.TextRange.Text = "3"
    }
    elseif ($sh.Id -eq 40) {
        # TextBox 39: "6" -> "2"
        $sh.TextFrame.TextRange.Text = "2"
    }
    elseif ($sh.Id -eq 50) {
        # Rectangle 49: "7" -> "1"
        $sh.TextFrame.TextRange.Text = "1"
    }
}
